$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.510.91"
$ws.Range("E2").Value = "  +0.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.914.30"
$ws.Range("E3").Value = "  -0.26%  "

$ws.Range("E4").Value = "  +0.67%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.10"
$ws.Range("E5").Value = "  -0.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  +0.60%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4828"
$ws.Range("E7").Value = "  +1.85%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4071"
$ws.Range("E8").Value = "  -0.61%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08136"
$ws.Range("E9").Value = "  +1.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.012"
$ws.Range("E10").Value = "  -0.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.39"
$ws.Range("E11").Value = "  +3.60%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.965.58"
$ws.Range("E12").Value = "  +1.49%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.000"
$ws.Range("E13").Value = "  +1.04%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.135"
$ws.Range("E14").Value = "  -0.44%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.23"
$ws.Range("E15").Value = "  +0.50%  "

$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.009"
$ws.Range("E16").Value = "  +0.82%  "

$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.06779"
$ws.Range("E17").Value = "  +2.81%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001039"
$ws.Range("E18").Value = "  +0.86%  "

$ws.Range("E19").Value = "  -0.44%  "

$ws.Range("E20").Value = "  +0.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.523.17"
$ws.Range("E21").Value = "  +0.08%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.626"
$ws.Range("E22").Value = "  +1.41%  "

$ws.Range("E23").Value = "  +2.36%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.185"
$ws.Range("E24").Value = "  -1.06%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.148.91"
$ws.Range("E25").Value = "  -0.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.42"
$ws.Range("E26").Value = "  +0.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.434"
$ws.Range("E27").Value = "  +7.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.06"
$ws.Range("E28").Value = "  +1.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.097"
$ws.Range("E29").Value = "  -1.67%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.71"
$ws.Range("E30").Value = "  +1.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.031"
$ws.Range("E31").Value = "  -3.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09564"
$ws.Range("E32").Value = "  +0.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.506"
$ws.Range("E33").Value = "  +1.66%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.567"
$ws.Range("E34").Value = "  +0.26%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.393"
$ws.Range("E35").Value = "  -3.08%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02268"
$ws.Range("E36").Value = "  +0.30%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06098"
$ws.Range("E37").Value = "  -0.27%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.176"
$ws.Range("E38").Value = "  +0.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5937"
$ws.Range("E39").Value = "  +0.44%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.76"
$ws.Range("E40").Value = "  +5.82%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.975"
$ws.Range("E41").Value = "  -4.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1859"
$ws.Range("E42").Value = "  +0.59%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.463"
$ws.Range("E43").Value = "  -3.85%  "

$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.282"
$ws.Range("E44").Value = "  -0.38%  "

$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.07696"
$ws.Range("E45").Value = "  -3.98%  "

$ws.Range("E46").Value = "  +1.90%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5575"
$ws.Range("E47").Value = "  +0.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.944"
$ws.Range("E48").Value = "  +0.42%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "115.61"
$ws.Range("E49").Value = "  +2.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "72.73"
$ws.Range("E50").Value = "  +1.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.053"
$ws.Range("E51").Value = "  +1.82%  "
